$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.29
$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 2.6
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 1.62
$ws.Range("T4").Value = 2.2
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 10
$ws.Range("AC4").Value = 5.5
$ws.Range("AR4").Value = 81
$ws.Range("AT4").Value = 2.2
$ws.Range("AZ4").Value = 151
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("P6").Value = 4.55
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.94
$ws.Range("H7").Value = 2.77
$ws.Range("I7").Value = 2.52
$ws.Range("K7").Value = 1.98
$ws.Range("L7").Value = 3.05
$ws.Range("N7").Value = 6.85
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 2.55
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.52
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.82
$ws.Range("X7").Value = 15.5
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 5.4
$ws.Range("AF7").Value = 70
$ws.Range("AH7").Value = 6.9
$ws.Range("AI7").Value = 12
$ws.Range("AJ7").Value = 9.5
$ws.Range("AL7").Value = 24
$ws.Range("AM7").Value = 35
$ws.Range("AT7").Value = 2.47
$ws.Range("AX7").Value = 13.5
$ws.Range("AY7").Value = 19.5
$ws.Range("AZ7").Value = 60
$ws.Range("BA7").Value = 90
$ws.Range("BB7").Value = 250
$ws.Range("G8").Value = 1.76
$ws.Range("I8").Value = 4.75
$ws.Range("J8").Value = 2.4
$ws.Range("N8").Value = 9.5
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("W8").Value = 7
$ws.Range("AC8").Value = 9.5
$ws.Range("AL8").Value = 41
$ws.Range("AU8").Value = 8.5
$ws.Range("AY8").Value = 34
$ws.Range("BB8").Value = 251
$ws.Range("G9").Value = 1.46
$ws.Range("BD9").Value = 176
$ws.Range("G10").Value = 2.2
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 3
$ws.Range("L10").Value = 4
$ws.Range("N10").Value = 8
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.62
$ws.Range("W10").Value = 6.5
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 9.5
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 21
$ws.Range("AC10").Value = 8
$ws.Range("AG10").Value = 351
$ws.Range("AH10").Value = 8.5
$ws.Range("AI10").Value = 15
$ws.Range("AJ10").Value = 12
$ws.Range("AK10").Value = 34
$ws.Range("AL10").Value = 29
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 13
$ws.Range("AP10").Value = 26
$ws.Range("AQ10").Value = 41
$ws.Range("AR10").Value = 67
$ws.Range("AS10").Value = 201
$ws.Range("AW10").Value = 5
$ws.Range("AX10").Value = 19
$ws.Range("AZ10").Value = 67
$ws.Range("BA10").Value = 101
$ws.Range("G12").Value = 2.05
$ws.Range("I12").Value = 3.6
$ws.Range("L12").Value = 4.5
$ws.Range("N12").Value = 7.5
$ws.Range("X12").Value = 8.5
$ws.Range("BA12").Value = 126
$ws.Range("G13").Value = 1.27
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 11
$ws.Range("L13").Value = 9
$ws.Range("AB13").Value = 29
$ws.Range("AF13").Value = 67
$ws.Range("AO13").Value = 6
$ws.Range("AW13").Value = 10
$ws.Range("AX13").Value = 41
$ws.Range("AY13").Value = 41
